# Updates cryptocurrency price (D) and 1h volume change (E) values on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" cells hold values that look like plain numbers (e.g. "7.00", "0.514").
# Force text format first so Excel does not auto-convert them to numbers and strip
# formatting/precision (matching the source data, which stores them as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values.
$ws.Range("D2").Value = "57.798.60"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "2.949.20"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("D5").Value = "555.34"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "133.28"
$ws.Range("E6").Value = "  +10.16%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +4.95%  "
$ws.Range("D9").Value = "2.944.84"
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("D11").Value = "4.83"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  +4.75%  "
$ws.Range("E13").Value = "  +5.70%  "
$ws.Range("D14").Value = "32.82"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "3.434.73"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").Value = "7.00"
$ws.Range("E17").Value = "  +9.33%  "
$ws.Range("D18").Value = "2.946.00"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "57.794.02"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").Value = "417.63"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "13.38"
$ws.Range("E21").Value = "  +5.16%  "
$ws.Range("D22").Value = "0.700"
$ws.Range("E22").Value = "  +7.90%  "
$ws.Range("D23").Value = "13.47"
$ws.Range("E23").Value = "  +7.76%  "
$ws.Range("D25").Value = "79.46"
$ws.Range("E25").Value = "  +4.12%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "2.50"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  +7.26%  "
$ws.Range("E30").Value = "  +6.89%  "
$ws.Range("D31").Value = "25.53"
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "0.0967"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  +7.18%  "
$ws.Range("D35").Value = "0.949"
$ws.Range("E35").Value = "  +7.11%  "
$ws.Range("D36").Value = "2.06"
$ws.Range("E36").Value = "  +2.91%  "
$ws.Range("D37").Value = "0.0₃0703"
$ws.Range("E37").Value = "  +15.06%  "
$ws.Range("D38").Value = "8.89"
$ws.Range("E38").Value = "  +7.12%  "
$ws.Range("D39").Value = "48.18"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  +16.60%  "
$ws.Range("D41").Value = "384.76"
$ws.Range("E41").Value = "  +6.74%  "
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("D43").Value = "0.0348"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").Value = "2.710.79"
$ws.Range("E44").Value = "  +4.72%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "124.89"
$ws.Range("E46").Value = "  +6.18%  "
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").Value = "22.97"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("E51").Value = "  +4.08%  "
